$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

# Update the Wild Card round distribution parameters (re-simulated values)
$ws.Range("B2").Value = "JSU(-1.2745210742677864, 1.3187981934704989, -0.40869963565298656, 3.2399218300920936)"
$ws.Range("C2").Value = "NCT(3.1307627292558085, 1.9601991755217627, -0.013226557929430723, 4.1548198732961446)"
$ws.Range("D2").Value = "NCT(2.6193011121542167, 1.3276918997431393, -0.4504814284224421, 2.646560748250657)"
$ws.Range("E2").Value = "NIG(1.622593757682529, 1.284080671362033, 2.8549548587141733, 6.738428301438214)"
